# Word COM-interop edit script
#
# Summary of the change (per the supplied OOXML diff):
#  1. Every paragraph in the document gains run-level / paragraph-mark
#     formatting: Times New Roman font plus sz/szCs = 24 (12pt), in
#     addition to whatever italics formatting already existed.
#  2. The "<link ...>" paragraph's single big run of text is broken up
#     into several runs (one per attribute-ish chunk), with
#     <w:proofErr w:type="spellStart"/>/<w:proofErr w:type="spellEnd"/>
#     pairs bracketing the "rel", "css" and "href" tokens (as Word's
#     spell-checker would mark them), and the inline _GoBack bookmark
#     that used to sit in the middle of that paragraph is removed.
#  3. A _GoBack bookmark is (re)created at the very end of the document,
#     immediately after the "</html>" run in the last non-empty
#     paragraph.
#  4. The final, empty trailing paragraph picks up paragraph-mark
#     formatting (Times New Roman, sz/szCs = 24) too.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: apply Times New Roman / 12pt (sz=24 half-points) to every
# paragraph's range. This naturally emits <w:rFonts .../>, <w:sz/> and
# <w:szCs/> on both the run(s) and the paragraph mark (w:pPr/w:rPr)
# while leaving existing direct formatting (the <w:i/><w:iCs/> pair)
# untouched, and without disturbing paragraph identity (w14:paraId /
# rsid attributes survive).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
    $p.Range.Font.Size = 12
    $p.Range.Font.SizeBi = 12
}

# ---------------------------------------------------------------------
# Step 2: rebuild the "<link ...>" paragraph (paragraph 5) as a series
# of runs, with proofErr spell-check markers around "rel", "css" and
# "href", and drop the bookmark that used to interrupt it.
# ---------------------------------------------------------------------
$rPrXml = '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/>'

$linkPara = $d.Paragraphs(5)
$linkRange = $linkPara.Range

$linkInner = '<w:pPr><w:rPr>' + $rPrXml + '</w:rPr></w:pPr>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:tab/><w:t xml:space="preserve">&lt;link </w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellStart"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>rel</w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellEnd"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>= “stylesheet” type= “text/</w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellStart"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>css</w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellEnd"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t xml:space="preserve">” </w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellStart"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>href</w:t></w:r>'
$linkInner += '<w:proofErr w:type="spellEnd"/>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>= “app.css</w:t></w:r>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>”</w:t></w:r>'
$linkInner += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>&gt;</w:t></w:r>'

$linkXml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $linkInner + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$linkRange.InsertXML($linkXml)

# ---------------------------------------------------------------------
# Step 3 + 4: rebuild the last two paragraphs together ("</html>" plus
# the trailing empty paragraph) so we can both append the _GoBack
# bookmark after "</html>" and give the empty trailing paragraph its
# own paragraph-mark formatting (InsertXML can't usefully target the
# lone end-of-document paragraph mark on its own).
# ---------------------------------------------------------------------
$htmlParaIndex = $d.Paragraphs.Count - 1
$htmlPara = $d.Paragraphs($htmlParaIndex)
$tailStart = $htmlPara.Range.Start
$tailEnd = $d.Content.End
$tailRange = $d.Range($tailStart, $tailEnd)

$tailInner1 = '<w:pPr><w:rPr>' + $rPrXml + '</w:rPr></w:pPr>'
$tailInner1 += '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>&lt;/html&gt;</w:t></w:r>'
$tailInner1 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$emptyRPrXml = '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/>'
$tailInner2 = '<w:pPr><w:rPr>' + $emptyRPrXml + '</w:rPr></w:pPr>'

$tailXml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $tailInner1 + '</w:p><w:p>' + $tailInner2 + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($tailXml)

Write-Output "edit complete"
